$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each of the 4 test cases has its title row ("Test Case..." in column B)
# and now gets a new "Passed" result cell added in column C.
$rows = @(1, 8, 17, 25)

# Highlight the existing test-case title cells with the light blue theme fill.
foreach ($row in $rows) {
    $ws.Cells.Item($row, 2).Interior.ThemeColor = 5
}

# Add the new "Passed" result cells with green font.
foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = "Passed"
    $cell.Font.Color = 5287936
}

# Update page setup (paper size / orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to C25.
$ws.Range("C25").Select() | Out-Null
